$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-12-15 03:54:50"
$wsDeDe.Range("H4").Value = "2016-12-15 03:54:50"
$wsZhCn.Range("H4").Value = "2016-12-15 03:54:36"
$wsZhCn.Range("L4").Value = "2016-12-15 03:55:31"
$wsDeDe.Range("L4").Value = "2016-12-15 03:55:50"
